# Update energy consumption values (rows 2-46, columns B and C)
# and remove the trailing rows 47-48 (performance improvement via priority queue
# changed the simulated values and shortened the series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value2 = 1.175807664859074
$ws.Cells.Item(2,3).Value2 = 1.13398605744893
$ws.Cells.Item(3,2).Value2 = 1.890136193836282
$ws.Cells.Item(3,3).Value2 = 2.372325209979794
$ws.Cells.Item(4,2).Value2 = 10.55753202155323
$ws.Cells.Item(4,3).Value2 = 3.446377943335892
$ws.Cells.Item(5,2).Value2 = 11.79587434895072
$ws.Cells.Item(5,3).Value2 = 4.91237614707967
$ws.Cells.Item(6,2).Value2 = 14.15087881023633
$ws.Cells.Item(6,3).Value2 = 5.952383050883469
$ws.Cells.Item(7,2).Value2 = 14.82028241398831
$ws.Cells.Item(7,3).Value2 = 7.195219367424396
$ws.Cells.Item(8,2).Value2 = 20.47054950121574
$ws.Cells.Item(8,3).Value2 = 8.835744457780587
$ws.Cells.Item(9,2).Value2 = 20.5788552329452
$ws.Cells.Item(9,3).Value2 = 10.02454788898532
$ws.Cells.Item(10,2).Value2 = 25.07479397661175
$ws.Cells.Item(10,3).Value2 = 11.29067658140405
$ws.Cells.Item(11,2).Value2 = 29.27589467309789
$ws.Cells.Item(11,3).Value2 = 12.70779535262184
$ws.Cells.Item(12,2).Value2 = 29.41898760410663
$ws.Cells.Item(12,3).Value2 = 13.93300158041861
$ws.Cells.Item(13,2).Value2 = 33.49016702753072
$ws.Cells.Item(13,3).Value2 = 15.11773759569614
$ws.Cells.Item(14,2).Value2 = 33.85918721095751
$ws.Cells.Item(14,3).Value2 = 16.34023581721498
$ws.Cells.Item(15,2).Value2 = 35.15338654512304
$ws.Cells.Item(15,3).Value2 = 17.7277679083827
$ws.Cells.Item(16,2).Value2 = 41.24098687204323
$ws.Cells.Item(16,3).Value2 = 18.83749431095422
$ws.Cells.Item(17,2).Value2 = 41.44217153928044
$ws.Cells.Item(17,3).Value2 = 20.03409627590943
$ws.Cells.Item(18,2).Value2 = 48.07760326123481
$ws.Cells.Item(18,3).Value2 = 21.52304310774542
$ws.Cells.Item(19,2).Value2 = 53.01559644521714
$ws.Cells.Item(19,3).Value2 = 22.67318962747508
$ws.Cells.Item(20,2).Value2 = 53.99210123626273
$ws.Cells.Item(20,3).Value2 = 23.83170405555349
$ws.Cells.Item(21,2).Value2 = 55.13018427940021
$ws.Cells.Item(21,3).Value2 = 25.00666805568065
$ws.Cells.Item(22,2).Value2 = 55.86198697690834
$ws.Cells.Item(22,3).Value2 = 26.11492771047858
$ws.Cells.Item(23,2).Value2 = 60.23622803199274
$ws.Cells.Item(23,3).Value2 = 27.09286354420058
$ws.Cells.Item(24,2).Value2 = 60.35363175490529
$ws.Cells.Item(24,3).Value2 = 28.28126436450325
$ws.Cells.Item(25,2).Value2 = 62.38020449461661
$ws.Cells.Item(25,3).Value2 = 29.61829324705363
$ws.Cells.Item(26,2).Value2 = 62.84669979531517
$ws.Cells.Item(26,3).Value2 = 30.8614845716239
$ws.Cells.Item(27,2).Value2 = 65.31098763305725
$ws.Cells.Item(27,3).Value2 = 32.68827877993377
$ws.Cells.Item(28,2).Value2 = 67.36485895618766
$ws.Cells.Item(28,3).Value2 = 33.83002510780443
$ws.Cells.Item(29,2).Value2 = 68.81517751687038
$ws.Cells.Item(29,3).Value2 = 34.9599537442171
$ws.Cells.Item(30,2).Value2 = 70.23275180143851
$ws.Cells.Item(30,3).Value2 = 36.07184543039403
$ws.Cells.Item(31,2).Value2 = 73.16266590261644
$ws.Cells.Item(31,3).Value2 = 37.22344491475403
$ws.Cells.Item(32,2).Value2 = 73.48527287084863
$ws.Cells.Item(32,3).Value2 = 38.41200214299177
$ws.Cells.Item(33,2).Value2 = 75.31741024587748
$ws.Cells.Item(33,3).Value2 = 40.09215433832844
$ws.Cells.Item(34,2).Value2 = 75.83721838253209
$ws.Cells.Item(34,3).Value2 = 41.47683125873904
$ws.Cells.Item(35,2).Value2 = 76.41803389444897
$ws.Cells.Item(35,3).Value2 = 42.69403660125447
$ws.Cells.Item(36,2).Value2 = 76.65015904251213
$ws.Cells.Item(36,3).Value2 = 44.25894942602997
$ws.Cells.Item(37,2).Value2 = 76.89922382769485
$ws.Cells.Item(37,3).Value2 = 45.60789493978878
$ws.Cells.Item(38,2).Value2 = 77.0401958641299
$ws.Cells.Item(38,3).Value2 = 46.97306082960549
$ws.Cells.Item(39,2).Value2 = 78.69519984507095
$ws.Cells.Item(39,3).Value2 = 48.19666393475044
$ws.Cells.Item(40,2).Value2 = 82.80540523957582
$ws.Cells.Item(40,3).Value2 = 49.4218742261416
$ws.Cells.Item(41,2).Value2 = 82.89968142660743
$ws.Cells.Item(41,3).Value2 = 50.65038361887588
$ws.Cells.Item(42,2).Value2 = 86.07006188353581
$ws.Cells.Item(42,3).Value2 = 51.91274043883428
$ws.Cells.Item(43,2).Value2 = 86.81130724751085
$ws.Cells.Item(43,3).Value2 = 53.11445583696225
$ws.Cells.Item(44,2).Value2 = 90.64465285293386
$ws.Cells.Item(44,3).Value2 = 54.31175463880344
$ws.Cells.Item(45,2).Value2 = 96.17760844488181
$ws.Cells.Item(45,3).Value2 = 55.49200645455882
$ws.Cells.Item(46,2).Value2 = 97.0685021886399
$ws.Cells.Item(46,3).Value2 = 56.73106741776657
# Remove the now-unused trailing rows (previously rows 47 and 48)
$ws.Rows("47:48").Delete()

Write-Output "Done updating values and removing rows 47-48"
